$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: Date, AMZN, RIOT, TSLA, SPCE, PLTR, GME, AMC, BBBY, BB, SNAP,
# NFLX, NVDA, AMD, SPOT, MSTR, SNOW, SQ, ROKU, ZM, SHOP, GOOG
$tickers = @("Date","AMZN","RIOT","TSLA","SPCE","PLTR","GME","AMC","BBBY","BB","SNAP","NFLX","NVDA","AMD","SPOT","MSTR","SNOW","SQ","ROKU","ZM","SHOP","GOOG")

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $tickers[$i]
}

$ws.Range("M12").Select()
